$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(80.40036015459947, 0.08040036015459946, 0.0008040036015459946)
    3 = @(119.5996398454005, 0.08040036015459946, 0.0008040036015459946)
    4 = @(80.40036015459947, 0.1195996398454005, 0.0008040036015459946)
    5 = @(119.5996398454005, 0.1195996398454005, 0.0008040036015459946)
    6 = @(80.40036015459947, 0.08040036015459946, 0.001195996398454005)
    7 = @(119.5996398454005, 0.08040036015459946, 0.001195996398454005)
    8 = @(80.40036015459947, 0.1195996398454005, 0.001195996398454005)
    9 = @(119.5996398454005, 0.1195996398454005, 0.001195996398454005)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
}
